$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value2 = "Elo"
$ws.Range("C2").Value2 = 1684
$ws.Range("C3").Value2 = 1678
$ws.Range("C4").Value2 = 1670
$ws.Range("B5").Value2 = "Nelson Bakerman"
$ws.Range("C5").Value2 = 1640
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.8
$ws.Range("B6").Value2 = "Chris Widgren"
$ws.Range("C6").Value2 = 1632
$ws.Range("D6").Value2 = 2
$ws.Range("E6").Value2 = 0
$ws.Range("F6").Value2 = 1
$ws.Range("B7").Value2 = "Eric Papa"
$ws.Range("C7").Value2 = 1631
$ws.Range("B8").Value2 = "Jeff Ziev"
$ws.Range("C8").Value2 = 1630
$ws.Range("D8").Value2 = 4
$ws.Range("E8").Value2 = 2
$ws.Range("F8").Value2 = 0.667
$ws.Range("B9").Value2 = "Ryan Leggette"
$ws.Range("C9").Value2 = 1627
$ws.Range("D9").Value2 = 3
$ws.Range("F9").Value2 = 0.75
$ws.Range("B10").Value2 = "Bob Sauchelli"
$ws.Range("A11").Value2 = 9
$ws.Range("B11").Value2 = "Paul Assad"
$ws.Range("C11").Value2 = 1617
$ws.Range("D11").Value2 = 2
$ws.Range("E11").Value2 = 1
$ws.Range("F11").Value2 = 0.667
$ws.Range("A12").Value2 = 11
$ws.Range("B12").Value2 = "David Chester"
$ws.Range("A13").Value2 = 11
$ws.Range("B14").Value2 = "Justin Goodfellow"
$ws.Range("B15").Value2 = "Anthony Buccellato"
$ws.Range("C15").Value2 = 1613
$ws.Range("D15").Value2 = 2
$ws.Range("F15").Value2 = 0.667
$ws.Range("A16").Value2 = 15
$ws.Range("B16").Value2 = "Damir Uzunic"
$ws.Range("C16").Value2 = 1602
$ws.Range("A17").Value2 = 16
$ws.Range("B17").Value2 = "Isaac Dunn"
$ws.Range("A18").Value2 = 16
$ws.Range("B18").Value2 = "Paul Jones"
$ws.Range("C18").Value2 = 1601
$ws.Range("D18").Value2 = 1
$ws.Range("E18").Value2 = 1
$ws.Range("A19").Value2 = 18
$ws.Range("B19").Value2 = "Kurowska"
$ws.Range("A20").Value2 = 18
$ws.Range("B20").Value2 = "Tom Witteman"
$ws.Range("A21").Value2 = 18
$ws.Range("B21").Value2 = "Mike Brady"
$ws.Range("B22").Value2 = "Ian Ainley"
$ws.Range("D22").Value2 = 2
$ws.Range("E22").Value2 = 2
$ws.Range("B23").Value2 = "Colin Kelly"
$ws.Range("A24").Value2 = 23
$ws.Range("B24").Value2 = "Adam Fratino"
$ws.Range("C24").Value2 = 1598
$ws.Range("D24").Value2 = 2
$ws.Range("E24").Value2 = 2
$ws.Range("B25").Value2 = "Judy O'Brien"
$ws.Range("C25").Value2 = 1597
$ws.Range("D25").Value2 = 1
$ws.Range("E25").Value2 = 1
$ws.Range("A26").Value2 = 25
$ws.Range("B26").Value2 = "Jeff Weber"
$ws.Range("C26").Value2 = 1588
$ws.Range("D26").Value2 = 2
$ws.Range("E26").Value2 = 3
$ws.Range("F26").Value2 = 0.4
$ws.Range("A27").Value2 = 26
$ws.Range("B27").Value2 = "George Brown"
$ws.Range("C27").Value2 = 1585
$ws.Range("E27").Value2 = 2
$ws.Range("F27").Value2 = 0.333
$ws.Range("A28").Value2 = 26
$ws.Range("B28").Value2 = "Steve Olsen"
$ws.Range("C28").Value2 = 1585
$ws.Range("D28").Value2 = 1
$ws.Range("E28").Value2 = 2
$ws.Range("F28").Value2 = 0.333
$ws.Range("B29").Value2 = "Dave Hitchings"
$ws.Range("D29").Value2 = 0
$ws.Range("E29").Value2 = 1
$ws.Range("F29").Value2 = 0
$ws.Range("A30").Value2 = 28
$ws.Range("A31").Value2 = 28
$ws.Range("B31").Value2 = "Chris Greene"
$ws.Range("B33").Value2 = "Colin Hanson"
$ws.Range("B34").Value2 = "Kofi Wilson"
$ws.Range("B36").Value2 = "Demelo"
$ws.Range("B37").Value2 = "Phil O'Brien"
$ws.Range("B38").Value2 = "Jeff Behrens"
$ws.Range("A39").Value2 = 38
$ws.Range("B39").Value2 = "Ben Cole"
$ws.Range("C39").Value2 = 1583
$ws.Range("D39").Value2 = 1
$ws.Range("E39").Value2 = 2
$ws.Range("F39").Value2 = 0.333
$ws.Range("A40").Value2 = 38
$ws.Range("B40").Value2 = "Roger Gibian"
$ws.Range("C40").Value2 = 1583
$ws.Range("D40").Value2 = 1
$ws.Range("E40").Value2 = 2
$ws.Range("F40").Value2 = 0.333
$ws.Range("A41").Value2 = 40
$ws.Range("C41").Value2 = 1583
$ws.Range("A42").Value2 = 40
$ws.Range("B42").Value2 = "Deb Czeresko"
$ws.Range("C42").Value2 = 1583
$ws.Range("A43").Value2 = 40
$ws.Range("B43").Value2 = "Amelia Burger"
$ws.Range("C43").Value2 = 1583
$ws.Range("B44").Value2 = "Matt Bird"
$ws.Range("C44").Value2 = 1582
$ws.Range("D44").Value2 = 1
$ws.Range("E44").Value2 = 2
$ws.Range("F44").Value2 = 0.333
$ws.Range("C45").Value2 = 1581
$ws.Range("A46").Value2 = 45
$ws.Range("B46").Value2 = "Pat Murphy"
$ws.Range("C46").Value2 = 1568
$ws.Range("D46").Value2 = 0
$ws.Range("F46").Value2 = 0
$ws.Range("A47").Value2 = 45
$ws.Range("B47").Value2 = "Alyssa Bird"
$ws.Range("A48").Value2 = 47
$ws.Range("B48").Value2 = "Dana Vandagriff"
$ws.Range("C48").Value2 = 1567
$ws.Range("A49").Value2 = 47
$ws.Range("C51").Value2 = 1553
